# TC04_Canine_Filter_FileType-RNASeqFile.xlsx
# "Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite"
#
# The CasesTab row's Neo4j query (B2) previously pulled an extra `Cohort`
# column via an OPTIONAL MATCH that isn't needed for this query - drop that
# trailing RETURN column (and the now-dangling trailing comma on the
# `Response to Treatment` line above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
   WHERE f.file_type IN ["RNA Sequence File"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Cells.Item(2, 2).Value2 = $newCasesQuery

# Row heights settle slightly differently once the CasesTab query loses a
# wrapped line (and the other two rows re-flow against the new font
# metrics too).
$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 259.2

# Column widths settle to very slightly narrower best-fit values too.
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 86.83
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.33
$ws.Columns.Item(5).ColumnWidth = 27.67

# Selection/active cell moves from B4 to B2 and the view scrolls back to
# the top of the sheet.
[void]$ws.Range("B2").Select()
